# Updated cryptos list — refresh Price (column D) and Volume(1h) (column E)
# values from the latest scrape, and correct the Cronos/Mantle row order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the new literal text that belongs in the given cell, exactly
# as scraped (including the two-space padding around percentages and the
# "thousand dot" formatting some prices use, e.g. "26.604.51").
$updates = @(
    @{ Cell = "D2"; Value = "26.604.51" },
    @{ Cell = "E2"; Value = "  -7.24%  " },
    @{ Cell = "D3"; Value = "1.695.25" },
    @{ Cell = "E3"; Value = "  -5.68%  " },
    @{ Cell = "E4"; Value = "  +0.30%  " },
    @{ Cell = "D5"; Value = "219.80" },
    @{ Cell = "E5"; Value = "  -5.16%  " },
    @{ Cell = "D6"; Value = "0.5099" },
    @{ Cell = "E6"; Value = "  -13.42%  " },
    @{ Cell = "E7"; Value = "  +0.22%  " },
    @{ Cell = "D8"; Value = "0.2654" },
    @{ Cell = "E8"; Value = "  -4.19%  " },
    @{ Cell = "D9"; Value = "22.16" },
    @{ Cell = "E9"; Value = "  -4.64%  " },
    @{ Cell = "D10"; Value = "0.06307" },
    @{ Cell = "E10"; Value = "  -6.87%  " },
    @{ Cell = "D11"; Value = "0.07367" },
    @{ Cell = "E11"; Value = "  -2.12%  " },
    @{ Cell = "D12"; Value = "1.696.08" },
    @{ Cell = "E12"; Value = "  -5.73%  " },
    @{ Cell = "D13"; Value = "4.526" },
    @{ Cell = "E13"; Value = "  -5.53%  " },
    @{ Cell = "D14"; Value = "0.5785" },
    @{ Cell = "E14"; Value = "  -6.42%  " },
    @{ Cell = "D15"; Value = "1.925.62" },
    @{ Cell = "E15"; Value = "  -5.64%  " },
    @{ Cell = "D16"; Value = "0.000008467" },
    @{ Cell = "E16"; Value = "  -7.04%  " },
    @{ Cell = "D17"; Value = "65.51" },
    @{ Cell = "E17"; Value = "  -13.03%  " },
    @{ Cell = "D18"; Value = "26.623.07" },
    @{ Cell = "E18"; Value = "  -7.09%  " },
    @{ Cell = "D19"; Value = "4.990" },
    @{ Cell = "E19"; Value = "  -8.84%  " },
    @{ Cell = "E20"; Value = "  +0.23%  " },
    @{ Cell = "E21"; Value = "  -4.54%  " },
    @{ Cell = "D22"; Value = "186.21" },
    @{ Cell = "E22"; Value = "  -11.59%  " },
    @{ Cell = "D23"; Value = "6.262" },
    @{ Cell = "E23"; Value = "  -8.27%  " },
    @{ Cell = "E24"; Value = "  +0.23%  " },
    @{ Cell = "D25"; Value = "144.66" },
    @{ Cell = "E25"; Value = "  -5.84%  " },
    @{ Cell = "D26"; Value = "7.476" },
    @{ Cell = "E26"; Value = "  -7.84%  " },
    @{ Cell = "D27"; Value = "0.1161" },
    @{ Cell = "E27"; Value = "  -8.09%  " },
    @{ Cell = "D28"; Value = "15.91" },
    @{ Cell = "E28"; Value = "  -3.31%  " },
    @{ Cell = "D29"; Value = "1.337" },
    @{ Cell = "E29"; Value = "  -6.06%  " },
    @{ Cell = "D30"; Value = "0.05712" },
    @{ Cell = "E30"; Value = "  -6.67%  " },
    @{ Cell = "D31"; Value = "1.342" },
    @{ Cell = "E31"; Value = "  -5.73%  " },
    @{ Cell = "D32"; Value = "3.519" },
    @{ Cell = "E32"; Value = "  -7.24%  " },
    @{ Cell = "D33"; Value = "3.503" },
    @{ Cell = "E33"; Value = "  -7.98%  " },
    @{ Cell = "D34"; Value = "1.643" },
    @{ Cell = "E34"; Value = "  -5.26%  " },
    @{ Cell = "E35"; Value = "  -2.87%  " },
    @{ Cell = "D36"; Value = "0.5994" },
    @{ Cell = "E36"; Value = "  -6.71%  " },
    @{ Cell = "D37"; Value = "2.362" },
    @{ Cell = "E37"; Value = "  -5.50%  " },
    @{ Cell = "D38"; Value = "2.688" },
    @{ Cell = "E38"; Value = "  -0.85%  " },
    @{ Cell = "D39"; Value = "0.01621" },
    @{ Cell = "E39"; Value = "  -4.49%  " },
    @{ Cell = "D40"; Value = "1.104.22" },
    @{ Cell = "E40"; Value = "  -3.38%  " },
    @{ Cell = "D41"; Value = "0.8587" },
    @{ Cell = "E41"; Value = "  -2.60%  " },
    @{ Cell = "D42"; Value = "5.835" },
    @{ Cell = "E42"; Value = "  -10.13%  " },
    @{ Cell = "E43"; Value = "  -0.24%  " },
    @{ Cell = "D44"; Value = "99.36" },
    @{ Cell = "E44"; Value = "  -0.82%  " },
    @{ Cell = "D45"; Value = "1.851.92" },
    @{ Cell = "E45"; Value = "  -5.10%  " },
    @{ Cell = "D46"; Value = "0.00000000117" },
    @{ Cell = "E46"; Value = "  +6.28%  " },
    @{ Cell = "D47"; Value = "56.62" },
    @{ Cell = "E47"; Value = "  -5.66%  " },
    @{ Cell = "D48"; Value = "1.004" },
    @{ Cell = "E48"; Value = "  +0.61%  " },
    @{ Cell = "D49"; Value = "8.116" },
    @{ Cell = "E49"; Value = "  -3.26%  " },
    @{ Cell = "B50"; Value = "Mantle" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" },
    @{ Cell = "D50"; Value = "0.4326" },
    @{ Cell = "E50"; Value = "  -3.43%  " },
    @{ Cell = "B51"; Value = "Cronos" },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" },
    @{ Cell = "D51"; Value = "0.05235" },
    @{ Cell = "E51"; Value = "  -4.60%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $val = $u.Value

    # Columns B/C/E are never ambiguous (coin names, URLs, "  +x.xx%  "
    # strings) so a plain assignment keeps them as text. Column D prices
    # are sometimes plain decimal numbers (e.g. "219.80", "0.00000000117")
    # which Excel would otherwise auto-convert to a Double on assignment,
    # silently dropping trailing zeros / switching to scientific notation.
    # Force those to stay text, then drop back to the default "Normal"
    # style so no stray number-format/quote-prefix styling is left behind.
    if ($u.Cell -match '^D\d+$' -and $val.Trim() -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
